# Atualização de bases das ligas, do dia: 01-06-2024 às 01:16
# Swap the full data (columns B through AD) between two rows, leaving
# column A (the running index) untouched, for the row pairs that were
# re-ordered in the source feed: (73,74) and (112,113).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

Swap-RowData 73 74 2 30
Swap-RowData 112 113 2 30
